$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the same format as C2 to D2 (matches style index 2 / numFmtId 22), then set value
$ws.Range("C2").Copy()
$ws.Range("D2").PasteSpecial(-4122)
$ws.Range("D2").Value = -1

# Update D3 and D4 values to -1
$ws.Range("D3").Value = -1
$ws.Range("D4").Value = -1

# Update the active selection to D2
$ws.Range("D2").Select()
